$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple value fixes -----------------------------------------------
# Fix the ConceptScheme URI (remove hyphen in m4m-20 -> m4m20)
$ws.Range("B1").Value = "http://purl.org/m4m20/subjects"

# Fix the PREFIX base URI (remove hyphen in m4m-20 -> m4m20)
$ws.Range("C3").Value = "http://purl.org/m4m20/subjects/"

# Update the "dct:modified" timestamp for the vocabulary
$ws.Range("B20").Value = "2022-06-04T21:44:07+00:00"

# --- Row 22 (column headers / mapping identifiers) ---------------------
# A new "qudt:unit" column (old C22) was removed, the remaining mapping
# columns shifted one place to the left, and two new date columns
# (dct:modified / dct:created) were inserted before dct:creator /
# dct:contributor, which moved right by one column (U,V -> V,W).
$ws.Range("C22").Value = 'skos:altLabel(separator=",")'
$ws.Range("D22").Value = "skos:definition@en"
$ws.Range("E22").Value = 'dct:source(separator=",")'
$ws.Range("F22").Value = 'skos:broader(lookupColumn="skos:prefLabel" separator=",")'
$ws.Range("G22").Value = 'skos:exactMatch(separator=",")'
$ws.Range("H22").Value = 'skos:closeMatch(separator=",")'
$ws.Range("I22").Value = "iop:hasProperty"
$ws.Range("J22").Value = "iop:hasObjectOfInterest"
$ws.Range("K22").Value = "iop:hasMatrix"
$ws.Range("L22").Value = 'iop:hasContextObject(separator=",")'
$ws.Range("M22").Value = 'iop:hasConstraint(separator=",")'
$ws.Range("N22").Value = 'puv:statistic(separator=",")'
$ws.Range("O22").Value = 'puv:usesMethod(separator=",")'
$ws.Range("P22").Value = 'sosa:madeBySensor(separator=",")'
$ws.Range("Q22").Value = 'puv:uom(separator=",")'
$ws.Range("R22").Value = "owl:deprecated^^xsd:boolean"
$ws.Range("S22").Value = "skos:editorialNote@en"
$ws.Range("T22").Value = "dct:modified^^xsd:date"
$ws.Range("U22").Value = "dct:created^^xsd:date"
$ws.Range("V22").Value = 'dct:creator(separator=",")'
$ws.Range("W22").Value = 'dct:contributor(separator=",")'

# --- Row 23 (vocab:1000 term) ------------------------------------------
# New test data was added for the prefLabel/altLabel/definition columns,
# and the old ORCID value (previously in U23, which is now a different
# mapped column) was cleared out.
$ws.Range("B23").Value = "test subject"
$ws.Range("C23").Value = "alt test subject"
$ws.Range("D23").Value = "some definition"
$ws.Range("U23").Value = ""

# --- Rows 24-27: clear the stale ORCID creator values -------------------
$ws.Range("U24").Value = ""
$ws.Range("U25").Value = ""
$ws.Range("U26").Value = ""
$ws.Range("U27").Value = ""
